$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.611.87"
$ws.Range("E2").Value = "  -2.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.005.59"
$ws.Range("E3").Value = "  -3.97%  "
$ws.Range("E4").Value = "  +1.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.87"
$ws.Range("E5").Value = "  -3.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.012"
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5021"
$ws.Range("E7").Value = "  -3.96%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4231"
$ws.Range("E8").Value = "  -3.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.15"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09031"
$ws.Range("E10").Value = "  -3.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.121"
$ws.Range("E11").Value = "  -3.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.33"
$ws.Range("E12").Value = "  -5.78%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.028.07"
$ws.Range("E13").Value = "  -2.80%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.064"
$ws.Range("E14").Value = "  -5.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.504"
$ws.Range("E15").Value = "  -5.58%  "
$ws.Range("E16").Value = "  +0.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.45"
$ws.Range("E17").Value = "  -6.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001115"
$ws.Range("E18").Value = "  -3.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06687"
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.67"
$ws.Range("E20").Value = "  -6.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.010"
$ws.Range("E21").Value = "  +0.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.966"
$ws.Range("E22").Value = "  -5.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.631.90"
$ws.Range("E23").Value = "  -2.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.04"
$ws.Range("E24").Value = "  -3.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.302"
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.01"
$ws.Range("E26").Value = "  -1.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.75"
$ws.Range("E27").Value = "  -4.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.403"
$ws.Range("E28").Value = "  -3.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.303"
$ws.Range("E29").Value = "  -8.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.43"
$ws.Range("E30").Value = "  -3.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.061"
$ws.Range("E31").Value = "  -6.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09980"
$ws.Range("E32").Value = "  -4.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.571"
$ws.Range("E33").Value = "  -5.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.853"
$ws.Range("E34").Value = "  -5.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.800"
$ws.Range("E35").Value = "  -1.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02480"
$ws.Range("E36").Value = "  -5.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.323"
$ws.Range("E37").Value = "  -8.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06400"
$ws.Range("E39").Value = "  -6.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6574"
$ws.Range("E40").Value = "  -5.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.73"
$ws.Range("E41").Value = "  -6.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2059"
$ws.Range("E42").Value = "  -6.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.010"
$ws.Range("E43").Value = "  +0.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6365"
$ws.Range("E44").Value = "  -6.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.48"
$ws.Range("E45").Value = "  -5.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.198"
$ws.Range("E46").Value = "  -5.37%  "
$ws.Range("E47").Value = "  -4.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.514"
$ws.Range("E48").Value = "  -3.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000334"
$ws.Range("E49").Value = "  -4.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07003"
$ws.Range("E50").Value = "  -3.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.130"
$ws.Range("E51").Value = "  -6.54%  "
